$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell's value while forcing it to remain Text (so purely
# numeric-looking strings like "581.38" aren't silently converted to a
# floating point Number by Excel's smart-entry). The NumberFormat is
# switched to Text ("@") just long enough for the assignment to "stick"
# as a string, then ClearFormats() strips the now-unneeded explicit
# style back off the cell (it keeps its already-entered text value/type),
# so the saved cell ends up with no style attribute at all -- same as
# every other untouched data cell in these columns.
function Set-TextValue($addr, $val) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
    $ws.Range($addr).ClearFormats()
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "67.778.49"
$ws.Range("E2").Value = "  -0.16%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.316.79"
$ws.Range("E3").Value = "  -1.84%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.03%  "

# Row 5 - BNB
Set-TextValue "D5" "581.38"
$ws.Range("E5").Value = "  -1.97%  "

# Row 6 - Solana
Set-TextValue "D6" "173.94"
$ws.Range("E6").Value = "  -7.27%  "

# Row 7 - USDC
Set-TextValue "D7" "0.999"
$ws.Range("E7").Value = "  +0.04%  "

# Row 8 - XRP
$ws.Range("E8").Value = "  -2.19%  "

# Row 9 - LidoStakedEther
$ws.Range("D9").Value = "3.311.94"
$ws.Range("E9").Value = "  -1.79%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  -5.10%  "

# Row 11 - Cardano
$ws.Range("E11").Value = "  -2.83%  "

# Row 12 - Avalanche
$ws.Range("E12").Value = "  -4.90%  "

# Row 13 - ShibaInu
$ws.Range("E13").Value = "  -3.29%  "

# Row 14 - BitcoinCash
Set-TextValue "D14" "663.16"
$ws.Range("E14").Value = "  +3.69%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "3.852.34"
$ws.Range("E15").Value = "  -1.76%  "

# Row 16 - Polkadot
Set-TextValue "D16" "8.36"
$ws.Range("E16").Value = "  -3.22%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "67.808.45"
$ws.Range("E17").Value = "  -0.12%  "

# Row 19 - WrappedEther
$ws.Range("D19").Value = "3.319.56"
$ws.Range("E19").Value = "  -1.81%  "

# Row 20 - Chainlink
Set-TextValue "D20" "17.43"
$ws.Range("E20").Value = "  -3.67%  "

# Row 21 - Uniswap
$ws.Range("E21").Value = "  -2.58%  "

# Row 22 - Polygon
$ws.Range("E22").Value = "  -2.80%  "

# Row 23 - Toncoin
Set-TextValue "D23" "5.43"
$ws.Range("E23").Value = "  +5.92%  "

# Row 24 - InternetComputer(DFINITY)
Set-TextValue "D24" "16.85"
$ws.Range("E24").Value = "  -6.37%  "

# Row 25 - Litecoin
Set-TextValue "D25" "97.24"
$ws.Range("E25").Value = "  -2.81%  "

# Row 26 - PancakeSwap
Set-TextValue "D26" "3.83"
$ws.Range("E26").Value = "  -5.33%  "

# Row 27 - ImmutableX
$ws.Range("E27").Value = "  -7.19%  "

# Row 28 - RenderToken
$ws.Range("E28").Value = "  -6.13%  "

# Row 29 - EthereumClassic
Set-TextValue "D29" "33.39"
$ws.Range("E29").Value = "  +1.96%  "

# Row 30 - Filecoin
$ws.Range("E30").Value = "  -3.92%  "

# Row 31 - NEARProtocol
Set-TextValue "D31" "7.23"
$ws.Range("E31").Value = "  +4.09%  "

# Row 32 - Bittensor
Set-TextValue "D32" "584.85"
$ws.Range("E32").Value = "  -5.08%  "

# Row 33 - Cosmos
$ws.Range("E33").Value = "  -1.90%  "

# Row 34 - Hedera
$ws.Range("E34").Value = "  -2.69%  "

# Row 35 - Maker
$ws.Range("D35").Value = "3.715.44"
$ws.Range("E35").Value = "  -8.20%  "

# Row 36 - Dai
$ws.Range("E36").Value = "  -0.16%  "

# Row 37 - OKB
Set-TextValue "D37" "56.68"
$ws.Range("E37").Value = "  +0.63%  "

# Row 38 - dogwifhat
Set-TextValue "D38" "3.30"
$ws.Range("E38").Value = "  -14.79%  "

# Row 39 - Kaspa
$ws.Range("E39").Value = "  -0.77%  "

# Row 40 - was InjectiveProtocol, now Fetch.AI
$ws.Range("B40").Value = "Fetch.AI"
$ws.Range("C40").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue "D40" "2.62"
$ws.Range("E40").Value = "  -7.31%  "

# Row 41 - was Fetch.AI, now InjectiveProtocol
$ws.Range("B41").Value = "InjectiveProtocol"
$ws.Range("C41").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue "D41" "32.32"
$ws.Range("E41").Value = "  -4.59%  "

# Row 42 - Stacks
$ws.Range("E42").Value = "  -5.58%  "

# Row 43 - TheGraph
$ws.Range("E43").Value = "  -3.89%  "

# Row 44 - was PEPE, now ApeXProtocol
$ws.Range("B44").Value = "ApeXProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
Set-TextValue "D44" "3.28"
$ws.Range("E44").Value = "  -4.01%  "

# Row 45 - was ApeXProtocol, now PEPE
$ws.Range("B45").Value = "PEPE"
$ws.Range("C45").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D45").Value = "0.0₃0662"
$ws.Range("E45").Value = "  -6.09%  "

# Row 46 - VeChain
Set-TextValue "D46" "0.0406"
$ws.Range("E46").Value = "  -4.31%  "

# Row 47 - ThetaToken
Set-TextValue "D47" "2.59"
$ws.Range("E47").Value = "  -0.52%  "

# Row 48 - Stellar
$ws.Range("E48").Value = "  -2.14%  "

# Row 49 - FirstDigitalUSD
$ws.Range("E49").Value = "  -0.03%  "

# Row 50 - Mantle
$ws.Range("E50").Value = "  -3.57%  "

# Row 51 - Monero
Set-TextValue "D51" "127.42"
$ws.Range("E51").Value = "  -0.54%  "
